$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '23.239.95'
$ws.Range('E2').Value = '  +0.94%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.602.81'
$ws.Range('E3').Value = '  +0.03%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('E5').Value = '  -0.16%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '303.45'
$ws.Range('E6').Value = '  +0.72%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3773'
$ws.Range('E7').Value = '  -0.26%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '51.96'
$ws.Range('E8').Value = '  +5.76%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3634'
$ws.Range('E9').Value = '  +0.08%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.276'
$ws.Range('E10').Value = '  +1.29%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.001'
$ws.Range('E11').Value = '  -0.11%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.08130'
$ws.Range('E12').Value = '  +0.22%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '22.83'
$ws.Range('E13').Value = '  +0.21%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.610'
$ws.Range('E14').Value = '  +0.36%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.428'
$ws.Range('E15').Value = '  +0.25%  '
$ws.Range('E16').Value = '  +0.41%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.605.38'
$ws.Range('E17').Value = '  +0.58%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '94.00'
$ws.Range('E18').Value = '  +2.23%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06925'
$ws.Range('E19').Value = '  +0.75%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '18.16'
$ws.Range('E20').Value = '  -0.37%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.535'
$ws.Range('E21').Value = '  -0.48%  '
$ws.Range('E22').Value = '  -0.29%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '12.97'
$ws.Range('E23').Value = '  -1.47%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '23.249.10'
$ws.Range('E24').Value = '  +0.96%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.037'
$ws.Range('E25').Value = '  +8.81%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.382'
$ws.Range('E26').Value = '  +0.39%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '21.22'
$ws.Range('E27').Value = '  +0.74%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '150.18'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.259'
$ws.Range('E29').Value = '  +0.08%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '134.63'
$ws.Range('E30').Value = '  +0.67%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.397'
$ws.Range('E31').Value = '  +4.02%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.769'
$ws.Range('E32').Value = '  -1.17%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.781.25'
$ws.Range('E33').Value = '  -0.14%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.9634'
$ws.Range('E34').Value = '  +0.16%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.07517'
$ws.Range('E35').Value = '  -1.51%  '
$ws.Range('B36').Value = 'FraxShare'
$ws.Range('C36').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '10.35'
$ws.Range('E36').Value = '  -0.11%  '
$ws.Range('B37').Value = 'VeChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02745'
$ws.Range('E37').Value = '  +1.24%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.2540'
$ws.Range('E38').Value = '  +0.32%  '
$ws.Range('B39').Value = 'InternetComputer(DFINITY)'
$ws.Range('C39').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '6.128'
$ws.Range('E39').Value = '  -2.55%  '
$ws.Range('B40').Value = 'Stellar'
$ws.Range('C40').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.08833'
$ws.Range('E40').Value = '  -0.15%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.396'
$ws.Range('E41').Value = '  +2.25%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.7117'
$ws.Range('E42').Value = '  +0.99%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '12.52'
$ws.Range('E43').Value = '  -0.12%  '
$ws.Range('E44').Value = '  +2.80%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.6544'
$ws.Range('E45').Value = '  -0.99%  '
$ws.Range('E46').Value = '  +0.26%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.9995'
$ws.Range('E47').Value = '  -0.12%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '4.015'
$ws.Range('E48').Value = '  +0.50%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '132.83'
$ws.Range('E49').Value = '  +0.19%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.07952'
$ws.Range('E50').Value = '  +0.48%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.206'
$ws.Range('E51').Value = '  -2.13%  '
